# Apply the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Tue Nov  7 19:38:55 UTC 2023 with GitHub Actions".
#
# Strategy:
#  - Text cells (Coin name / Link) are assigned directly as strings.
#  - Percentage cells (column E) keep their leading/trailing double-space
#    padding, so Excel never mistakes them for numbers and they stay text.
#  - Price cells (column D) sometimes look like genuine numbers (e.g.
#    "246.49"), which Excel's COM automatism would silently coerce into a
#    floating point value. To keep them as plain text (matching the
#    original "35.408.39"-style values), each D-cell is briefly switched
#    to the "@" (Text) number format before the value is written, then
#    ClearFormats() removes that temporary formatting again so the cell's
#    style stays the same as before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '35.409.46'
$c.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  +0.25%  '

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.894.04'
$c.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  -0.95%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.72%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '246.49'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -3.25%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.692'
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -4.38%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.77%  '

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '43.78'
$c.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +7.47%  '

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.350'
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -5.33%  '

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.0738'
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -3.29%  '

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.0971'
$c.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -1.59%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '13.12'
$c.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +2.54%  '

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '2.166.82'
$c.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -0.95%  '

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.722'
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -0.34%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '4.90'
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -1.14%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '1.898.87'
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.95%  '

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '35.373.86'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +0.26%  '

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '73.46'
$c.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -1.38%  '

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0823'
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -3.87%  '

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '246.54'
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +1.02%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '12.80'
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -2.35%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '4.95'
$c.ClearFormats()

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.81%  '

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +6.22%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -11.90%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '165.26'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -1.15%  '

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '8.46'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -2.35%  '

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '18.34'
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -2.41%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -4.43%  '

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '4.128.41'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -0.03%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +7.82%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -2.94%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '0.0579'
$c.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -1.46%  '

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +0.36%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.852'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -6.67%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -2.36%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -19.80%  '

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '17.20'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.47%  '

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.0680'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +5.24%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '97.29'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +0.36%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -3.32%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -2.98%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -2.87%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '1.285.32'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -3.87%  '

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.0801'
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +7.28%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.01%  '

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -0.76%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +1.84%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -5.53%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '43.08'
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -4.87%  '
